$wb = $excel.ActiveWorkbook

# --- Sheet "Data": append the new weekly observation row (row 94) ---
$wsData = $wb.Worksheets.Item("Data")

# Copy format of the last existing data row (93) down into the new row (94),
# then overwrite the values with the new observation.
$wsData.Range("A93:B93").Copy($wsData.Range("A94:B94"))
$wsData.Range("A94").Value2 = 45119
$wsData.Range("B94").Value2 = 8296.923000000001

# --- Sheet "SeriesInfo": refresh the metadata pulled from FRED ---
$wsInfo = $wb.Worksheets.Item("SeriesInfo")

function Set-TextValue($range, $text) {
    # Force the cell to be treated as literal text so Excel does not
    # reinterpret date-like strings (e.g. "2023-07-20") as date serials.
    $range.NumberFormat = "@"
    $range.Value2 = $text
    $range.Style = "Normal"
}

Set-TextValue $wsInfo.Range("B3") "2023-07-20"
Set-TextValue $wsInfo.Range("B4") "2023-07-20"
Set-TextValue $wsInfo.Range("B7") "2023-07-12"
Set-TextValue $wsInfo.Range("B14") "2023-07-13 15:33:35-05"
